# Martan_1.xlsx edit script
# - drop the two extra (empty) sheets, keep a single sheet renamed "Martan_1"
# - touch up several dialogue lines in column B
# - apply word-wrap + a wider column B + matching row heights + a period font
# - move the active selection to B7

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- remove "Лист2" / "Лист3", keep only the first sheet ---
while ($wb.Worksheets.Count -gt 1) {
    [void]$wb.Worksheets.Item($wb.Worksheets.Count).Delete()
}
$ws.Name = "Martan_1"

# --- corrected / re-edited dialogue text ---
$ws.Cells.Item(1, 2).Value = "*Мартан Килмови, твой непосредственный начальник, смотрит на тебя. Вполне возможно, что он скучает. А может быть и нет. Трудно сказать. У него очень усталые глаза. От жилетки и рубашки веет сильным запахом типичного мужского одеколона.* Капитан, рад вас видеть. Докладывайте."
$ws.Cells.Item(4, 2).Value = "Мартан… сэр. Кстати, насчет этого. Что мне можно говорить? Чего говорить нельзя?"
$ws.Cells.Item(5, 2).Value = "*Килмови качает головой* Это зависит от того, что вам удалось узнать. Точно могу сказать только то, что про «Аврору-1» лучше помалкивать в любом случае. А так… черт возьми... это ты у нас герой мирового масштаба, я-то что могу тебе посоветовать?"
$ws.Cells.Item(7, 2).Value = "Увидим. Тебе пора идти. Удачи и... да. Твое место за общим столом — все еще твое."
$ws.Cells.Item(10, 2).Value = "Ну, хоть что-то хорошее…"

# --- formatting: wrap column B, widen it, set the period-correct font ---
$ws.Columns.Item(2).ColumnWidth = 66.140625
$ws.Range("A1:C10").WrapText = $true
$ws.Cells.Font.Name = "Arial Cyr"
$ws.Cells.Font.Size = 10

$ws.Rows.Item(1).RowHeight = 63.75
$ws.Rows.Item(2).RowHeight = 25.5
$ws.Rows.Item(3).RowHeight = 25.5
$ws.Rows.Item(4).RowHeight = 25.5
$ws.Rows.Item(5).RowHeight = 51
$ws.Rows.Item(7).RowHeight = 25.5

# --- legacy (Excel 2003-style) page margins ---
$ws.PageSetup.LeftMargin = 0.75 * 72
$ws.PageSetup.RightMargin = 0.75 * 72
$ws.PageSetup.TopMargin = 1 * 72
$ws.PageSetup.BottomMargin = 1 * 72
$ws.PageSetup.HeaderMargin = 0.5 * 72
$ws.PageSetup.FooterMargin = 0.5 * 72

# --- selection / view ---
[void]$ws.Range("B7").Select()

Write-Output "done"
